# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
# with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.882.82"
$ws.Range("D3").Value = "2.598.17"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'551.77"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").Value = "'143.34"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +5.70%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "3.054.02"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "58.848.05"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "'20.89"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "2.607.34"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "'4.47"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "'337.46"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "'10.06"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").Value = "'6.17"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'66.94"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").Value = "'7.15"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "0.0₃0754"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "'5.99"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("D32").Value = "'154.73"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("D33").Value = "'18.96"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "'3.93"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("D35").Value = "'0.891"
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("D36").Value = "'1.13"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "'36.98"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").Value = "'0.829"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'284.00"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'0.599"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").Value = "'0.0959"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "'0.0534"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "'0.0227"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "1.945.39"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "'118.40"
$ws.Range("E49").Value = "  +6.43%  "
$ws.Range("D50").Value = "'17.89"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").Value = "'4.42"
$ws.Range("E51").Value = "  -3.15%  "
